$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value2 = '42.895.60'
$ws.Range('E2').Value2 = '  -1.14%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value2 = '2.310.57'
$ws.Range('E3').Value2 = '  +0.01%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value2 = '1.00'
$ws.Range('E4').Value2 = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '303.05'
$ws.Range('E5').Value2 = '  -1.84%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value2 = '100.15'
$ws.Range('E6').Value2 = '  -4.76%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value2 = '0.505'
$ws.Range('E7').Value2 = '  -3.79%  '
$ws.Range('E8').Value2 = '  +0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '0.504'
$ws.Range('E9').Value2 = '  -2.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value2 = '34.87'
$ws.Range('E10').Value2 = '  -2.75%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value2 = '0.0792'
$ws.Range('E11').Value2 = '  -2.25%  '
$ws.Range('E12').Value2 = '  +0.48%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '6.73'
$ws.Range('E13').Value2 = '  -3.41%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value2 = '2.669.55'
$ws.Range('E14').Value2 = '  -0.14%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '15.67'
$ws.Range('E15').Value2 = '  +3.57%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value2 = '2.326.35'
$ws.Range('E16').Value2 = '  +0.53%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value2 = '0.799'
$ws.Range('E17').Value2 = '  -0.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '42.808.43'
$ws.Range('E18').Value2 = '  -1.25%  '
$ws.Range('B19').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value2 = '11.69'
$ws.Range('E19').Value2 = '  -1.52%  '
$ws.Range('B20').Value2 = 'ShibaInu'
$ws.Range('C20').Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value2 = '0.0₃0906'
$ws.Range('E20').Value2 = '  -1.72%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '6.06'
$ws.Range('E21').Value2 = '  -2.25%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value2 = '67.94'
$ws.Range('E22').Value2 = '  -0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '235.79'
$ws.Range('E23').Value2 = '  -1.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value2 = '1.96'
$ws.Range('E24').Value2 = '  -3.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '2.51'
$ws.Range('E25').Value2 = '  -3.12%  '
$ws.Range('E26').Value2 = '  -0.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '24.86'
$ws.Range('E27').Value2 = '  +0.01%  '
$ws.Range('E28').Value2 = '  -1.43%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value2 = '34.54'
$ws.Range('E29').Value2 = '  -4.62%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value2 = '165.74'
$ws.Range('E30').Value2 = '  +1.95%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '9.14'
$ws.Range('E31').Value2 = '  -4.69%  '
$ws.Range('E32').Value2 = '  -0.06%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value2 = '5.02'
$ws.Range('E33').Value2 = '  -4.08%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value2 = '2.41'
$ws.Range('E34').Value2 = '  -4.85%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value2 = '4.46'
$ws.Range('E35').Value2 = '  -2.74%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '16.83'
$ws.Range('E36').Value2 = '  -8.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '0.0698'
$ws.Range('E37').Value2 = '  -4.59%  '
$ws.Range('E38').Value2 = '  -3.64%  '
$ws.Range('E39').Value2 = '  -2.99%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value2 = '0.101'
$ws.Range('E40').Value2 = '  -4.66%  '
$ws.Range('E41').Value2 = '  -3.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '2.54'
$ws.Range('E42').Value2 = '  +1.73%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '1.971.38'
$ws.Range('E43').Value2 = '  +0.28%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '0.0281'
$ws.Range('E44').Value2 = '  -2.83%  '
$ws.Range('B45').Value2 = 'EnergySwap'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '18.50'
$ws.Range('E45').Value2 = '  -1.50%  '
$ws.Range('B46').Value2 = 'FraxShare'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '10.27'
$ws.Range('E46').Value2 = '  -0.31%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '2.88'
$ws.Range('E47').Value2 = '  -5.86%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '55.66'
$ws.Range('E48').Value2 = '  -4.01%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value2 = '2.533.37'
$ws.Range('E49').Value2 = '  -0.16%  '
$ws.Range('E50').Value2 = '  -3.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '4.68'
$ws.Range('E51').Value2 = '  +0.38%  '
